# Apply edits to PALMARES PRODUIT worksheet per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ---
# (ColumnWidth is quantized internally to 1/6-character steps by this
# COM layer, so the inputs below are chosen to land on the raw OOXML
# widths 12 and ~15.6 as closely as that quantization allows.)
$ws.Columns.Item(3).ColumnWidth = 11.166666666666666
$ws.Columns.Item(4).ColumnWidth = 14.766666666666666

# --- Header row text ---
$ws.Range("C1").Value = "QUANTITE"
$ws.Range("D1").Value = "MONTANT TTC"

# --- Prepare formatting for the new rows before writing data ---
# Row 13 currently holds the TOTAUX formatting (bold / filled style).
# Copy that formatting to the new TOTAUX row (16) first.
$ws.Range("A13:D13").Copy()
$ws.Range("A16:D16").PasteSpecial(-4122)

# Copy the standard data-row formatting (from row 2) onto row 13
# (which becomes a normal data row) and onto the two brand-new rows 14-15.
$ws.Range("A2:D2").Copy()
$ws.Range("A13:D13").PasteSpecial(-4122)
$ws.Range("A14:D14").PasteSpecial(-4122)
$ws.Range("A15:D15").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Write the final data set (rows 2-15) ---
$ws.Range("A2").Value = "BAF008"
$ws.Range("B2").Value = "RIBEYE / ENTRECOTE"
$ws.Range("C2").Value = 3.134
$ws.Range("D2").Value = 173968.34

$ws.Range("A3").Value = "BLPCG001"
$ws.Range("B3").Value = "BLANC DE POULET CONGELE"
$ws.Range("C3").Value = 3.4
$ws.Range("D3").Value = 98600

$ws.Range("A4").Value = "MERL0001"
$ws.Range("B4").Value = " LAMB MERGUEZ / MERGUEZ D'AGNEAU "
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 60880

$ws.Range("A5").Value = "BBQS0001"
$ws.Range("B5").Value = "BBQ SAUSAGE / SAUCISSE BBQ"
$ws.Range("C5").Value = 0.824
$ws.Range("D5").Value = 54911.36

$ws.Range("A6").Value = "LAML0001"
$ws.Range("B6").Value = "LAMB LEG CHOPS / GIGOT TRANCHE"
$ws.Range("C6").Value = 0.688
$ws.Range("D6").Value = 42621.6

$ws.Range("A7").Value = "BEESK002"
$ws.Range("B7").Value = "BEEF SKEWERS - WITH FAT"
$ws.Range("C7").Value = 0.768
$ws.Range("D7").Value = 38330.88

$ws.Range("A8").Value = "GOAR0002"
$ws.Range("B8").Value = "GOAT RIB / COTE DE CHEVRE"
$ws.Range("C8").Value = 0.674
$ws.Range("D8").Value = 37103.7

$ws.Range("A9").Value = "ALPCG001"
$ws.Range("B9").Value = "AILES DE POULET CONGELE"
$ws.Range("C9").Value = 1.222
$ws.Range("D9").Value = 29939

$ws.Range("A10").Value = "PICB0001"
$ws.Range("B10").Value = "PICANHA BRAZILIAN CUT"
$ws.Range("C10").Value = 0.444
$ws.Range("D10").Value = 29254.05

$ws.Range("A11").Value = "MERC0001"
$ws.Range("B11").Value = "MERGUEZ CONGELE"
$ws.Range("C11").Value = 0.862
$ws.Range("D11").Value = 25860

$ws.Range("A12").Value = "BEESI001"
$ws.Range("B12").Value = "BEEF SIRLOIN / ALOYAU"
$ws.Range("C12").Value = 0.404
$ws.Range("D12").Value = 21294.84

$ws.Range("A13").Value = "BRIA0001"
$ws.Range("B13").Value = "AMERICAN BRISKET"
$ws.Range("C13").Value = 0.322
$ws.Range("D13").Value = 15790.88

$ws.Range("A14").Value = "FOICG001"
$ws.Range("B14").Value = "FOIE CONGELE"
$ws.Range("C14").Value = 0.292
$ws.Range("D14").Value = 5110

$ws.Range("A15").Value = "SACHET_0"
$ws.Range("B15").Value = "PLASTIC BAGS"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 500

# --- Write TOTAUX row (16) ---
$ws.Range("A16").Value = "TOTAUX"
$ws.Range("C16").Value = 15.034
$ws.Range("D16").Value = 634164.65
